# Insert a new bridge slide ("YOU HAVE THREE OPTIONS") as slide 3,
# between the problem-framing slide (THE GAP IS GETTING WORSE) and the
# solution slide (DESCRIBE THE MISSION. GET THE SILICON.)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Insert the new, blank slide at position 3 (ppLayoutBlank = 12).
#    Grab shape/style templates only AFTER this insertion so that slide
#    index references below are stable (inserting shifts later slides).
# ---------------------------------------------------------------------
$ns = $p.Slides.Add(3, 12)

# Background fill 2D2D2D
$ns.FollowMasterBackground = $false
$ns.Background.Fill.Visible = $true
$ns.Background.Fill.ForeColor.RGB = 0x2D2D2D

# Template shapes to copy the theme-styled "p:style" boilerplate from -
# after insertion, slide4 (orig slide3) & slide5 (orig slide4) hold them.
$barTemplate = $p.Slides.Item(4).Shapes.Item(2)          # "Rectangle 2"
$roundRectTemplate = $p.Slides.Item(6).Shapes.Item(3)     # "Rounded Rectangle 3"

# ---------------------------------------------------------------------
# 2. Shape 2: Title textbox - "YOU HAVE THREE OPTIONS"
# ---------------------------------------------------------------------
$title = $ns.Shapes.AddTextbox(1, 72, 57.6, 792, 72)
$title.Name = "TextBox 1"
$title.Fill.Visible = $false
$title.TextFrame.WordWrap = 1
$title.TextFrame.AutoSize = 1
$ttr = $title.TextFrame.TextRange
$ttr.Text = "YOU HAVE THREE OPTIONS"
$ttr.ParagraphFormat.Alignment = 1
$ttr.Font.Size = 44
$ttr.Font.Bold = $true
$ttr.Font.Name = "Calibri"
$ttr.Font.Color.RGB = 0xFFFFFF
$title.Height = 72

# ---------------------------------------------------------------------
# 3. Shape 3: decorative red accent bar ("Rectangle 2")
# ---------------------------------------------------------------------
$barTemplate.Copy()
$bar = $ns.Shapes.Paste().Item(1)
$bar.Name = "Rectangle 2"
$bar.Left = 72
$bar.Top = 144
$bar.Width = 144
$bar.Height = 3
$bar.Fill.ForeColor.RGB = 0x3D4DE8

# ---------------------------------------------------------------------
# 4. Shape 4: "Rounded Rectangle 3" - Option 01 - Wait for the roadmap
# ---------------------------------------------------------------------
$roundRectTemplate.Copy()
$rr1 = $ns.Shapes.Paste().Item(1)
$rr1.Name = "Rounded Rectangle 3"
$rr1.Left = 72
$rr1.Top = 187.2
$rr1.Width = 237.6
$rr1.Height = 201.6
$rr1.Fill.ForeColor.RGB = 0x3A3A3A
$rr1.TextFrame.WordWrap = 1
$rr1.TextFrame.MarginLeft = 18
$rr1.TextFrame.MarginRight = 18
$rr1.TextFrame.MarginTop = 18
$rr1.TextFrame.MarginBottom = 18

$t1 = $rr1.TextFrame.TextRange
$t1.Text = "01`rWait for the roadmap`rYou ship when they ship. If they ship."

$p1a = $t1.Paragraphs(1, 1)
$p1a.ParagraphFormat.Alignment = 1
$p1a.Font.Size = 28
$p1a.Font.Bold = $true
$p1a.Font.Name = "Calibri"
$p1a.Font.Color.RGB = 0x3D4DE8

$p1b = $t1.Paragraphs(2, 1)
$p1b.ParagraphFormat.Alignment = 1
$p1b.ParagraphFormat.SpaceBefore = 10
$p1b.Font.Size = 22
$p1b.Font.Bold = $true
$p1b.Font.Name = "Calibri"
$p1b.Font.Color.RGB = 0xFFFFFF

$p1c = $t1.Paragraphs(3, 1)
$p1c.ParagraphFormat.Alignment = 1
$p1c.ParagraphFormat.SpaceBefore = 12
$p1c.Font.Size = 18
$p1c.Font.Bold = $false
$p1c.Font.Name = "Calibri"
$p1c.Font.Color.RGB = 0xCCCCCC

# ---------------------------------------------------------------------
# 5. Shape 5: "Rounded Rectangle 4" - Option 02 - Throw money at it
# ---------------------------------------------------------------------
$roundRectTemplate.Copy()
$rr2 = $ns.Shapes.Paste().Item(1)
$rr2.Name = "Rounded Rectangle 4"
$rr2.Left = 345.6
$rr2.Top = 187.2
$rr2.Width = 237.6
$rr2.Height = 201.6
$rr2.Fill.ForeColor.RGB = 0x3A3A3A
$rr2.TextFrame.WordWrap = 1
$rr2.TextFrame.MarginLeft = 18
$rr2.TextFrame.MarginRight = 18
$rr2.TextFrame.MarginTop = 18
$rr2.TextFrame.MarginBottom = 18

$t2 = $rr2.TextFrame.TextRange
$t2.Text = "02`rThrow money at it`rOverpay today, repeat next generation."

$p2a = $t2.Paragraphs(1, 1)
$p2a.ParagraphFormat.Alignment = 1
$p2a.Font.Size = 28
$p2a.Font.Bold = $true
$p2a.Font.Name = "Calibri"
$p2a.Font.Color.RGB = 0x3D4DE8

$p2b = $t2.Paragraphs(2, 1)
$p2b.ParagraphFormat.Alignment = 1
$p2b.ParagraphFormat.SpaceBefore = 10
$p2b.Font.Size = 22
$p2b.Font.Bold = $true
$p2b.Font.Name = "Calibri"
$p2b.Font.Color.RGB = 0xFFFFFF

$p2c = $t2.Paragraphs(3, 1)
$p2c.ParagraphFormat.Alignment = 1
$p2c.ParagraphFormat.SpaceBefore = 12
$p2c.Font.Size = 18
$p2c.Font.Bold = $false
$p2c.Font.Name = "Calibri"
$p2c.Font.Color.RGB = 0xCCCCCC

# ---------------------------------------------------------------------
# 6. Shape 6: "Rounded Rectangle 5" - Option 03 - Design your own compute
# ---------------------------------------------------------------------
$roundRectTemplate.Copy()
$rr3 = $ns.Shapes.Paste().Item(1)
$rr3.Name = "Rounded Rectangle 5"
$rr3.Left = 619.2
$rr3.Top = 187.2
$rr3.Width = 237.6
$rr3.Height = 201.6
$rr3.Fill.ForeColor.RGB = 0x2A4A2A
$rr3.TextFrame.WordWrap = 1
$rr3.TextFrame.MarginLeft = 18
$rr3.TextFrame.MarginRight = 18
$rr3.TextFrame.MarginTop = 18
$rr3.TextFrame.MarginBottom = 18

$t3 = $rr3.TextFrame.TextRange
$t3.Text = "03`rDesign your own compute`rMatch silicon to mission, own the trajectory."

$p3a = $t3.Paragraphs(1, 1)
$p3a.ParagraphFormat.Alignment = 1
$p3a.Font.Size = 28
$p3a.Font.Bold = $true
$p3a.Font.Name = "Calibri"
$p3a.Font.Color.RGB = 0xE99B3D

$p3b = $t3.Paragraphs(2, 1)
$p3b.ParagraphFormat.Alignment = 1
$p3b.ParagraphFormat.SpaceBefore = 10
$p3b.Font.Size = 22
$p3b.Font.Bold = $true
$p3b.Font.Name = "Calibri"
$p3b.Font.Color.RGB = 0xFFFFFF

$p3c = $t3.Paragraphs(3, 1)
$p3c.ParagraphFormat.Alignment = 1
$p3c.ParagraphFormat.SpaceBefore = 12
$p3c.Font.Size = 18
$p3c.Font.Bold = $false
$p3c.Font.Name = "Calibri"
$p3c.Font.Color.RGB = 0xCCCCCC

# ---------------------------------------------------------------------
# 7. Shape 7: Bottom caption textbox ("TextBox 6")
# ---------------------------------------------------------------------
$caption = $ns.Shapes.AddTextbox(1, 72, 424.8, 792, 72)
$caption.Name = "TextBox 6"
$caption.Fill.Visible = $false
$caption.TextFrame.WordWrap = 1
$caption.TextFrame.AutoSize = 1

$ctr = $caption.TextFrame.TextRange
$ctr.Text = "Option 3 used to take 50 engineers and 18 months.  "
$ctr.ParagraphFormat.Alignment = 1
$ctr.Font.Size = 22
$ctr.Font.Bold = $false
$ctr.Font.Name = "Calibri"
$ctr.Font.Color.RGB = 0x998888

$lenBefore = $ctr.Length
$full = $ctr.InsertAfter("It doesn't anymore.")
$run2 = $full.Characters($lenBefore + 1, 20)
$run2.Font.Size = 22
$run2.Font.Bold = $true
$run2.Font.Name = "Calibri"
$run2.Font.Color.RGB = 0xFFFFFF

$caption.Height = 72
